$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (column D) cells whose new values would otherwise
# be auto-parsed by Excel as numbers (losing exact text representation, e.g. trailing zeros).
$textCells = @("D5", "D6", "D8", "D9", "D13", "D17", "D18", "D19", "D21", "D25", "D27", "D30", "D31", "D37", "D38", "D39", "D40", "D44", "D45", "D46", "D49")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values (coin names/links swapped for rows 13/14 and 38/39,
# plus updated price and 1h-volume-change figures throughout).
$ws.Range('D2').Value = '34.419.27'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '1.790.62'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '226.05'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '0.556'
$ws.Range('E6').Value = '  +1.47%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '32.68'
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('D9').Value = '0.296'
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '2.049.12'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '11.07'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.783.80'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').Value = '34.390.41'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('D17').Value = '4.25'
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('D18').Value = '68.71'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('D19').Value = '246.65'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  +2.56%  '
$ws.Range('D21').Value = '11.20'
$ws.Range('E21').Value = '  +3.14%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').Value = '164.43'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').Value = '16.48'
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('E28').Value = '  +2.44%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = '1.23'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').Value = '3.79'
$ws.Range('E31').Value = '  +3.40%  '
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('E33').Value = '  +6.77%  '
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('D35').Value = '1.417.93'
$ws.Range('E35').Value = '  -2.08%  '
$ws.Range('E36').Value = '  +5.16%  '
$ws.Range('D37').Value = '0.670'
$ws.Range('E37').Value = '  +2.92%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0192'
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.06'
$ws.Range('E39').Value = '  +1.55%  '
$ws.Range('D40').Value = '84.62'
$ws.Range('E40').Value = '  +5.20%  '
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('E43').Value = '  +2.16%  '
$ws.Range('D44').Value = '13.66'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('D45').Value = '0.0524'
$ws.Range('E45').Value = '  +3.02%  '
$ws.Range('D46').Value = '6.04'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').Value = '1.949.30'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').Value = '105.39'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('E50').Value = '  -3.40%  '
$ws.Range('E51').Value = '  -0.03%  '
